$d = $word.ActiveDocument

# Locate the "Introduction" heading paragraph, then the empty paragraph right after it.
$rng = $d.Content
$found = $rng.Find.Execute("Introduction", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Introduction heading not found"
}
$headingStart = $rng.Start

$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $headingStart) {
        $headingIndex = $i
        break
    }
}
if ($headingIndex -eq -1) {
    throw "Could not resolve heading paragraph index"
}

$targetIndex = $headingIndex + 1
$p = $d.Paragraphs.Item($targetIndex)
$r = $d.Range($p.Range.Start, $p.Range.End)

# Replace the empty paragraph in place with one containing the justified intro
# text, three leading soft hyphens, and matching Times New Roman run formatting.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00AF1EAF" w:rsidRPr="00AF1EAF" w:rsidRDefault="00AF1EAF" w:rsidP="00AF1EAF"><w:pPr><w:spacing w:after="120"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:softHyphen/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:softHyphen/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:softHyphen/><w:t xml:space="preserve">There exist two standard methods of calculating the flight characteristics associated with a given aircraft: flight simulation of the whole aircraft through the use of potential flow theory and computational fluid dynamics and wind tunnel testing on a model representation of a given aircraft, where the flow characteristics are scaled up to the full size aircraft through the use of dimensional analysis. Each method has its advantages and disadvantages, with the </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr><w:t>former offering a quick, oftentimes highly accurate, representation all of the characteristics associated with an aircraft’s lift, drag, and moment coefficients for a wide range of angle of attacks. The accuracy of this method, though substantial and improving all the time, is highly reliant on the proper running of and collection of data from real-life wind tunnel testing. This report serves as documentation for one particular wind tunnel test that was run within the Low-Speed Tunnel within Oliver Hall on the campus of Saint Louis University.</w:t></w:r></w:p>'

$r.InsertXML($xml)

Write-Host "Updated paragraph" $targetIndex "now contains" $d.Paragraphs.Item($targetIndex).Range.Characters.Count "characters."
